$d = $word.ActiveDocument

# Locate the point in the "time management" paragraph right after
# "managemen" (i.e. just before the final "t" of "management") and
# insert a stray "s" there, turning "management" into "managemenst".
$r = $d.Content
[void]$r.Find.Execute("time managemen", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertPoint = $d.Range($r.End, $r.End)
$insertPoint.InsertBefore("s")

# Force the freshly inserted "s" to live in its own run (matching the
# target document) by round-tripping it through a temporary bookmark;
# adding/removing a bookmark at that one-character range splits the
# surrounding run without altering any visible formatting.
$sRange = $d.Range($r.End, $r.End + 1)
[void]$d.Bookmarks.Add("__tmp_split__", $sRange)
[void]$d.Bookmarks("__tmp_split__").Delete()
